$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "SubjectKind (URN, SubjectKind, Attribute, Value) : Statement;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SubjectKind (URN, SubjectReifiedKind, PredicateAttribute, ObjectValue) : Statement;",
    2)

$d.Content.Find.Execute(
    "PredicateKind (URN, PredicateKind, Attribute, Value) : Statement;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PredicateKind (URN, PredicateReifiedKind, SubjectAttribute, ObjectValue) : Statement;",
    2)

$d.Content.Find.Execute(
    "ObjectKind (URN, SubjectKind, Attribute, Value) : Statement;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ObjectKind (URN, ObjectReifiedKind, PredicateAttribute, SubjectValue) : Statement;",
    2)
